# Sample Project / Main.xlsx — "Rules" sheet
# The rule in row 11 (B11) used to read the shared string "R40"; the
# project was re-saved with that cell's value changed to "1" (still a
# text value, not a number) while every other aspect of the cell
# (its style/border, the other cells in the row, etc.) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new literal to be stored as text (the source cell always held
# a string), then write the new value.
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "1"
